# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps on the
# per-language report sheets to reflect the newly generated report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-24 07:05:12"
$wsZhCn.Range("G2").Value = "2016-02-24 07:06:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-24 07:05:33"
$wsDeDe.Range("G2").Value = "2016-02-24 07:06:48"
